# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps written during report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (row 2)
$wsOverview.Range("G2").Value = "2016-08-31 01:09:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn.Range("H2").Value = "2016-08-31 01:09:14"
$wsZhCn.Range("K2").Value = "2016-08-31 01:09:32"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsDeDe.Range("H2").Value = "2016-08-31 01:09:19"
$wsDeDe.Range("K2").Value = "2016-08-31 01:09:39"
